$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2026-02-08 -> 2026-02-09, i.e. 46061 -> 46062) for every data row (2..283).
$ws.Range("C2:C283").Value = 46062
